$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) values for the
# cryptos list, as refreshed by the scheduled GitHub Actions data pull.
# Numeric-looking price strings are forced to text so Excel's automatic
# type conversion doesn't turn them into floating point numbers (which
# would also silently drop significant trailing zeros).

$ws.Range('D2').Value = '43.715.73'
$ws.Range('E2').Value = '  -0.31%  '
$ws.Range('D3').Value = '2.312.66'
$ws.Range('E3').Value = '  +3.93%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '269.86'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.48%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '93.39'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +6.90%  '
$ws.Range('E7').Value = '  +0.86%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.621'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.45%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '44.92'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.08%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0937'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +1.09%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '8.06'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +7.01%  '
$ws.Range('D14').Value = '2.658.32'
$ws.Range('E14').Value = '  +3.86%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.32'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +3.43%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.852'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +7.24%  '
$ws.Range('D17').Value = '2.321.53'
$ws.Range('E17').Value = '  +4.70%  '
$ws.Range('D18').Value = '43.697.18'
$ws.Range('E18').Value = '  -0.26%  '
$ws.Range('E19').Value = '  +1.83%  '
$ws.Range('E20').Value = '  +3.90%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.36'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +1.85%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '239.64'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +3.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.28'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -4.48%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.69'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +8.51%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.50'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -5.08%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '11.26'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +3.97%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.33'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +5.64%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '3.37'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -4.68%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '38.96'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.28%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '22.55'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +9.50%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '171.64'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.92%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0894'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.39%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.59'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.75%  '
$ws.Range('E35').Value = '  +1.49%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.111'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.15%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.51'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.85%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0349'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.33%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.40'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.49%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.235'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +15.45%  '
$ws.Range('E41').Value = '  +8.74%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '12.18'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.74%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.30'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +15.01%  '
$ws.Range('E44').Value = '  +1.78%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '61.33'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -6.58%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '8.91'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +7.18%  '
$ws.Range('E47').Value = '  +2.99%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '100.16'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.13%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.20'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.57%  '
$ws.Range('D50').Value = '2.536.85'
$ws.Range('E50').Value = '  +3.83%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.427'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -3.16%  '
